# Generate Report for Handoff
#
# The localization-status report is regenerated: the "ae0d94ec" file, which
# had been handed back, is now ready for a new handoff (new handoff
# timestamps), while the "f0067453" file keeps its existing "handed back"
# status. Rows are re-sorted so ae0d94ec now appears in row 3 and f0067453
# in row 2 across all three sheets.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay {
    param($ws, [string]$addr, [string]$text)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value = "2016-32-12 08:32:52"

$wsOverview.Range("A3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-34-12 08:34:00"

Set-HyperlinkDisplay $wsOverview '$A$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
Set-HyperlinkDisplay $wsOverview '$A$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-12 08:32:28"
$wsZh.Range("F2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
$wsZh.Range("G2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-03-12 08:33:24"
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-12 08:33:56"
$wsZh.Range("F3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
$wsZh.Range("G3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-03-12 08:33:24"
$wsZh.Range("I3").Value = "Include"

Set-HyperlinkDisplay $wsZh '$A$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
Set-HyperlinkDisplay $wsZh '$D$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.zh-cn.xlf"
Set-HyperlinkDisplay $wsZh '$F$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
Set-HyperlinkDisplay $wsZh '$G$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.zh-cn.xlf"

Set-HyperlinkDisplay $wsZh '$A$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
Set-HyperlinkDisplay $wsZh '$D$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.zh-cn.xlf"
Set-HyperlinkDisplay $wsZh '$F$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
Set-HyperlinkDisplay $wsZh '$G$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-12 08:32:52"
$wsDe.Range("F2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
$wsDe.Range("G2").Value = "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.de-de.xlf"
$wsDe.Range("H2").Value = "2016-03-12 08:33:30"
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-12 08:34:00"
$wsDe.Range("F3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
$wsDe.Range("G3").Value = "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-03-12 08:33:30"
$wsDe.Range("I3").Value = "Include"

Set-HyperlinkDisplay $wsDe '$A$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
Set-HyperlinkDisplay $wsDe '$D$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.de-de.xlf"
Set-HyperlinkDisplay $wsDe '$F$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.md"
Set-HyperlinkDisplay $wsDe '$G$2' "f0067453-4126-4cb2-81ef-ade8d03f27c3.c10901518470bf49267ab6ef4bddc88faba38a6f.de-de.xlf"

Set-HyperlinkDisplay $wsDe '$A$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
Set-HyperlinkDisplay $wsDe '$D$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.de-de.xlf"
Set-HyperlinkDisplay $wsDe '$F$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.md"
Set-HyperlinkDisplay $wsDe '$G$3' "ae0d94ec-2646-4100-8858-101a6a503f0e.1e15a0e2b32aad7b1581d34b1d20e1277064190c.de-de.xlf"
